# Adapt column header formatting to respective input file names.
# Renames header suffixes "_old" -> "_FV2210" and "_new" -> "_FV2304",
# wraps the data range in an Excel Table, and freezes the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename the header cells in row 1 (A1:U1) ---
$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count
$lastCol = $usedRange.Columns.Count

for ($c = 1; $c -le $lastCol; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $val = $cell.Value()
    if ($val -ne $null) {
        $newVal = $val -replace '_old$', '_FV2210'
        $newVal = $newVal -replace '_new$', '_FV2304'
        if ($newVal -ne $val) {
            $cell.Value = $newVal
        }
    }
}

# --- 2. Wrap the data range into an Excel Table (ListObject) ---
$rng = $ws.Range($ws.Cells.Item(1, 1), $ws.Cells.Item($lastRow, $lastCol))
$tbl = $ws.ListObjects.Add(1, $rng, $null, 1)
$tbl.Name = "Table1"

# --- 3. Freeze the header row (pane split after row 1) ---
$ws.Activate()
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true

$wb.Save() | Out-Null
